# Update forest data - 2025-12-23 12:20
#
# The "New" worksheet holds freshly scraped listings that get promoted into
# "Previously added" once they age out, while brand-new listings land in
# "New". This edit moves the 4 existing rows on "New" (rows 2-5) down to the
# bottom of "Previously added", and replaces "New" with a single freshly
# scraped row.

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# Data that currently lives on "New" rows 2-5 (about to be archived into
# "Previously added" rows 344-347) together with the hyperlink target for
# each row's column A cell.
# ---------------------------------------------------------------------
$movedRows = @(
    @{ Url = "https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/vecumnieku-pag/lebhx.html"; Price = "9 800 €"; Region = "Bauska un raj."; Area = "3.32 ha."; Cadastre = "40940130127"; Date = 46013.58541666667 },
    @{ Url = "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/piedrujas-pag/kjohf.html"; Price = "44 000 €"; Region = "Krāslava un raj."; Area = "9.50 ha."; Cadastre = "60840050087, 162"; Date = 46010.92638888889 },
    @{ Url = "https://www.ss.com/msg/lv/real-estate/wood/tukums-and-reg/smardes-pag/dcegf.html"; Price = "2 200 €"; Region = "Tukums un raj."; Area = "0.08 ha."; Cadastre = "90820050176"; Date = 46012.92291666666 },
    @{ Url = "https://www.ss.com/msg/lv/real-estate/wood/tukums-and-reg/smardes-pag/kecid.html"; Price = "15 000 €"; Region = "Tukums un raj."; Area = "4 ha."; Cadastre = "90820050174"; Date = 46012.92152777778 }
)

# New single row that replaces everything on the "New" sheet.
$freshRow = @{ Url = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/pilskalnes-pag/lfcem.html"; Price = "7 500 €"; Region = "Daugavpils un raj."; Area = "17000 m²"; Cadastre = ""; Date = 46014.59444444445 }

# ---------------------------------------------------------------------
# 1) Append the 4 moved rows to the bottom of "Previously added"
#    (rows 344-347), re-using the formatting of the last existing row.
# ---------------------------------------------------------------------
$lastRow = $wsPrev.UsedRange.Rows.Count
$firstNewRow = $lastRow + 1

for ($i = 0; $i -lt $movedRows.Count; $i++) {
    $destRow = $firstNewRow + $i
    $srcRange = $wsPrev.Range("A" + $lastRow + ":F" + $lastRow)
    $dstRange = $wsPrev.Range("A" + $destRow + ":F" + $destRow)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)

    $data = $movedRows[$i]
    $wsPrev.Range("A" + $destRow).Value = $data.Url
    $wsPrev.Range("B" + $destRow).Value = $data.Price
    $wsPrev.Range("C" + $destRow).Value = $data.Region
    $wsPrev.Range("D" + $destRow).Value = $data.Area
    $wsPrev.Range("E" + $destRow).Value = $data.Cadastre
    $wsPrev.Range("F" + $destRow).Value = $data.Date

    $wsPrev.Hyperlinks.Add($wsPrev.Range("A" + $destRow), $data.Url)
}

# ---------------------------------------------------------------------
# 2) Clear "New": drop the old hyperlinks + rows 3-5, then overwrite row 2
#    with the single freshly scraped listing.
# ---------------------------------------------------------------------
$wsNew.Range("A2").Hyperlinks.Delete()
$wsNew.Range("A3").Hyperlinks.Delete()
$wsNew.Range("A4").Hyperlinks.Delete()
$wsNew.Range("A5").Hyperlinks.Delete()
$wsNew.Rows("3:5").Delete()

$wsNew.Range("A2").Value = $freshRow.Url
$wsNew.Range("B2").Value = $freshRow.Price
$wsNew.Range("C2").Value = $freshRow.Region
$wsNew.Range("D2").Value = $freshRow.Area
$wsNew.Range("E2").Value = $freshRow.Cadastre
$wsNew.Range("F2").Value = $freshRow.Date

$wsNew.Hyperlinks.Add($wsNew.Range("A2"), $freshRow.Url)

Write-Host "Forest data updated: moved 4 rows to 'Previously added', added 1 fresh row to 'New'."
